# Weekly data refresh for "Fruta, Vega Modelo de Temuco - Tuna":
# a new day's record is inserted at row 31 (pushing the existing
# rows 31:75 down to 32:76, so last row of data becomes row 76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 31, shifting rows 31-75 down
# to 32-76 (dimension grows from A1:T75 to A1:T76).
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new week's record.
$ws.Range("A31").Value2 = 10
$ws.Range("B31").Value2 = "Vega Modelo de Temuco"
$ws.Range("C31").Value2 = "La Araucanía"
$ws.Range("D31").Value2 = 44994
$ws.Range("E31").Value2 = 9
$ws.Range("F31").Value2 = "Fruta"
$ws.Range("G31").Value2 = 100107
$ws.Range("H31").Value2 = "Otros"
$ws.Range("I31").Value2 = 100107011
$ws.Range("J31").Value2 = "Tuna"
$ws.Range("K31").Value2 = "Sin especificar"
$ws.Range("L31").Value2 = "Primera"
$ws.Range("M31").Value2 = 100
$ws.Range("N31").Value2 = 20000
$ws.Range("O31").Value2 = 20000
$ws.Range("P31").Value2 = 20000
$ws.Range("Q31").Value2 = "`$/caja 16 kilos"
$ws.Range("R31").Value2 = "Provincia de Los Andes"
$ws.Range("S31").Value2 = 1250
$ws.Range("T31").Value2 = 16
